$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Utopía) corrections
$ws.Range("AD2").Value = -43.4597587435138
$ws.Range("AN2").Value = 1.15558321354578
$ws.Range("AV2").Value = 100.768467681994

# Row 3 (Distopía) corrections
$ws.Range("F3").Value = -37.4293802769236
$ws.Range("AV3").Value = -582.586667681994
$ws.Range("AZ3").Value = 12.078182486723
$ws.Range("BE3").Value = -0.54627143911055
